# Natmi following Dr Hou advice
# Recompute the Lama2-Itgb1 LR-pair table: cluster-pairs are now the full
# cross product of {ECs, FAPs, sCs} sending x {ECs, FAPs, sCs} target (9 rows)
# instead of only the 3 same-target-as-row-2/5 pairs (7 rows incl. header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs (Lama2/Itgb1)
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Lama2"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.471482
$ws.Cells.Item(2, 8).Value = 4.414446
$ws.Cells.Item(2, 9).Value = 0.004946458467382327
$ws.Cells.Item(2, 10).Value = 0.004946458467382326
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 117.044563
$ws.Cells.Item(2, 14).Value = 351.133689
$ws.Cells.Item(2, 15).Value = 0.3245365645427815
$ws.Cells.Item(2, 16).Value = 0.3245365645427815
$ws.Cells.Item(2, 17).Value = 172.228967652366
$ws.Cells.Item(2, 18).Value = 1550.060708871294
$ws.Cells.Item(2, 19).Value = 0.001605306637657813
$ws.Cells.Item(2, 20).Value = 0.001605306637657812

# Row 3: FAPs -> FAPs (Lama2/Itgb1)
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Lama2"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.471482
$ws.Cells.Item(3, 8).Value = 4.414446
$ws.Cells.Item(3, 9).Value = 0.004946458467382327
$ws.Cells.Item(3, 10).Value = 0.004946458467382326
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 101.5800373333333
$ws.Cells.Item(3, 14).Value = 304.740112
$ws.Cells.Item(3, 15).Value = 0.281657135515876
$ws.Cells.Item(3, 16).Value = 0.281657135515876
$ws.Cells.Item(3, 17).Value = 149.473196495328
$ws.Cells.Item(3, 18).Value = 1345.258768457952
$ws.Cells.Item(3, 19).Value = 0.001393205322871156
$ws.Cells.Item(3, 20).Value = 0.001393205322871156

# Row 4: FAPs -> sCs (Lama2/Itgb1)
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Lama2"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.471482
$ws.Cells.Item(4, 8).Value = 4.414446
$ws.Cells.Item(4, 9).Value = 0.004946458467382327
$ws.Cells.Item(4, 10).Value = 0.004946458467382326
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 142.0267893333333
$ws.Cells.Item(4, 14).Value = 426.080368
$ws.Cells.Item(4, 15).Value = 0.3938062999413425
$ws.Cells.Item(4, 16).Value = 0.3938062999413425
$ws.Cells.Item(4, 17).Value = 208.989864021792
$ws.Cells.Item(4, 18).Value = 1880.908776196128
$ws.Cells.Item(4, 19).Value = 0.001947946506853358
$ws.Cells.Item(4, 20).Value = 0.001947946506853358

# Row 5: ECs -> ECs (Lama2/Itgb1)
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Lama2"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 264.7713316666666
$ws.Cells.Item(5, 8).Value = 794.313995
$ws.Cells.Item(5, 9).Value = 0.8900417371348598
$ws.Cells.Item(5, 10).Value = 0.8900417371348596
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 117.044563
$ws.Cells.Item(5, 14).Value = 351.133689
$ws.Cells.Item(5, 15).Value = 0.3245365645427815
$ws.Cells.Item(5, 16).Value = 0.3245365645427815
$ws.Cells.Item(5, 17).Value = 30990.04480985306
$ws.Cells.Item(5, 18).Value = 278910.4032886776
$ws.Cells.Item(5, 19).Value = 0.2888510876694367
$ws.Cells.Item(5, 20).Value = 0.2888510876694367

# Row 6: ECs -> FAPs (Lama2/Itgb1)
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Lama2"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 264.7713316666666
$ws.Cells.Item(6, 8).Value = 794.313995
$ws.Cells.Item(6, 9).Value = 0.8900417371348598
$ws.Cells.Item(6, 10).Value = 0.8900417371348596
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 101.5800373333333
$ws.Cells.Item(6, 14).Value = 304.740112
$ws.Cells.Item(6, 15).Value = 0.281657135515876
$ws.Cells.Item(6, 16).Value = 0.281657135515876
$ws.Cells.Item(6, 17).Value = 26895.48175549638
$ws.Cells.Item(6, 18).Value = 242059.3357994674
$ws.Cells.Item(6, 19).Value = 0.2506866061709789
$ws.Cells.Item(6, 20).Value = 0.2506866061709788

# Row 7: ECs -> sCs (Lama2/Itgb1)
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Lama2"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 264.7713316666666
$ws.Cells.Item(7, 8).Value = 794.313995
$ws.Cells.Item(7, 9).Value = 0.8900417371348598
$ws.Cells.Item(7, 10).Value = 0.8900417371348596
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 142.0267893333333
$ws.Cells.Item(7, 14).Value = 426.080368
$ws.Cells.Item(7, 15).Value = 0.3938062999413425
$ws.Cells.Item(7, 16).Value = 0.3938062999413425
$ws.Cells.Item(7, 17).Value = 37604.6221441278
$ws.Cells.Item(7, 18).Value = 338441.5992971502
$ws.Cells.Item(7, 19).Value = 0.3505040432944441
$ws.Cells.Item(7, 20).Value = 0.3505040432944441

# Row 8: sCs -> ECs (Lama2/Itgb1)
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Lama2"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 31.239114
$ws.Cells.Item(8, 8).Value = 93.717342
$ws.Cells.Item(8, 9).Value = 0.105011804397758
$ws.Cells.Item(8, 10).Value = 0.105011804397758
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 117.044563
$ws.Cells.Item(8, 14).Value = 351.133689
$ws.Cells.Item(8, 15).Value = 0.3245365645427815
$ws.Cells.Item(8, 16).Value = 0.3245365645427815
$ws.Cells.Item(8, 17).Value = 3656.368446637182
$ws.Cells.Item(8, 18).Value = 32907.31601973464
$ws.Cells.Item(8, 19).Value = 0.03408017023568695
$ws.Cells.Item(8, 20).Value = 0.03408017023568694

# Row 9: sCs -> FAPs (Lama2/Itgb1)
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Lama2"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 31.239114
$ws.Cells.Item(9, 8).Value = 93.717342
$ws.Cells.Item(9, 9).Value = 0.105011804397758
$ws.Cells.Item(9, 10).Value = 0.105011804397758
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 101.5800373333333
$ws.Cells.Item(9, 14).Value = 304.740112
$ws.Cells.Item(9, 15).Value = 0.281657135515876
$ws.Cells.Item(9, 16).Value = 0.281657135515876
$ws.Cells.Item(9, 17).Value = 3173.270366380256
$ws.Cells.Item(9, 18).Value = 28559.43329742231
$ws.Cells.Item(9, 19).Value = 0.029577324022026
$ws.Cells.Item(9, 20).Value = 0.02957732402202599

# Row 10: sCs -> sCs (Lama2/Itgb1)
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Lama2"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 31.239114
$ws.Cells.Item(10, 8).Value = 93.717342
$ws.Cells.Item(10, 9).Value = 0.105011804397758
$ws.Cells.Item(10, 10).Value = 0.105011804397758
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 142.0267893333333
$ws.Cells.Item(10, 14).Value = 426.080368
$ws.Cells.Item(10, 15).Value = 0.3938062999413425
$ws.Cells.Item(10, 16).Value = 0.3938062999413425
$ws.Cells.Item(10, 17).Value = 4436.791063037985
$ws.Cells.Item(10, 18).Value = 39931.11956734186
$ws.Cells.Item(10, 19).Value = 0.0413543101400451
$ws.Cells.Item(10, 20).Value = 0.04135431014004508
